# Refresh the "cryptos" price table (Price / Volume(1h) columns, plus two
# rows whose ranking order flipped) to match the latest GitHub Actions run.
#
# Note: several Price-column values look like plain numbers (e.g. "0.578",
# "36.19"). Assigning those to Range.Value directly would make Excel parse
# them as numeric, changing both the stored type and the cell's number
# format. To keep them as literal text (matching the source data, which
# uses dotted thousand separators inconsistently, e.g. "42.787.09"), we
# briefly force NumberFormat to Text ("@") before the assignment and clear
# the format again afterwards so the cell's style reverts to its original
# (default/general) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.787.09"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.527.57"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("E4").Value = "  -0.08%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "314.71"
$r.ClearFormats()
$ws.Range("E5").Value = "  +1.93%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "95.72"
$r.ClearFormats()
$ws.Range("E6").Value = "  -0.18%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.578"
$r.ClearFormats()
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.97%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "36.19"
$r.ClearFormats()
$ws.Range("E10").Value = "  -0.69%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0810"
$r.ClearFormats()
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").Value = "2.915.02"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.508.24"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "15.28"
$r.ClearFormats()
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "42.846.70"
$ws.Range("E18").Value = "  +0.74%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "6.77"
$r.ClearFormats()
$ws.Range("E19").Value = "  +4.28%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "12.86"
$r.ClearFormats()
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  -0.91%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "69.81"
$r.ClearFormats()
$ws.Range("E22").Value = "  -2.30%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "252.70"
$r.ClearFormats()
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +4.14%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "40.83"
$r.ClearFormats()
$ws.Range("E29").Value = "  +8.60%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  +0.38%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "157.46"
$r.ClearFormats()
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  +4.30%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "19.39"
$r.ClearFormats()
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "3.32"
$r.ClearFormats()
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "2.68"
$r.ClearFormats()
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  -1.06%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "23.42"
$r.ClearFormats()
$ws.Range("E40").Value = "  -5.44%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "2.33"
$r.ClearFormats()
$ws.Range("E41").Value = "  +15.59%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("E44").Value = "  +0.23%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "3.32"
$r.ClearFormats()
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("D46").Value = "2.051.16"
$ws.Range("E46").Value = "  +1.21%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "84.92"
$r.ClearFormats()
$ws.Range("E47").Value = "  +0.74%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "107.43"
$r.ClearFormats()
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("E49").Value = "  -0.25%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "75.35"
$r.ClearFormats()
$ws.Range("E50").Value = "  +3.75%  "
$ws.Range("D51").Value = "2.769.18"
$ws.Range("E51").Value = "  +0.55%  "
